$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style (used in column B) from an existing row down to
# the new rows so the new date cells keep the same number format / style
# index as the rest of the table.
$ws.Range("B137").Copy($ws.Range("B138:B140"))

# Row 138 - Phil / Sauntering Hippo
$ws.Cells.Item(138, 1).Value = "Phil"
$ws.Cells.Item(138, 2).Value = 45473
$ws.Cells.Item(138, 3).Value = "Workout"
$ws.Cells.Item(138, 4).Value = 44
$ws.Cells.Item(138, 5).Value = 0
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 44
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = "Sauntering Hippo"
$ws.Cells.Item(138, 13).Value = 3

# Row 139 - Matt / Agile Antelope
$ws.Cells.Item(139, 1).Value = "Matt"
$ws.Cells.Item(139, 2).Value = 45474
$ws.Cells.Item(139, 3).Value = "Workout"
$ws.Cells.Item(139, 4).Value = 31
$ws.Cells.Item(139, 5).Value = 0
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 23
$ws.Cells.Item(139, 8).Value = 8
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = "Agile Antelope"
$ws.Cells.Item(139, 13).Value = 4

# Row 140 - Jeremiah / Agile Antelope
$ws.Cells.Item(140, 1).Value = "Jeremiah"
$ws.Cells.Item(140, 2).Value = 45474
$ws.Cells.Item(140, 3).Value = "Workout"
$ws.Cells.Item(140, 4).Value = 47
$ws.Cells.Item(140, 5).Value = 0
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 40
$ws.Cells.Item(140, 8).Value = 6
$ws.Cells.Item(140, 9).Value = 1
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = "Agile Antelope"
$ws.Cells.Item(140, 13).Value = 4

# Scroll / select, mirroring what Excel records after appending rows and
# moving to the next blank row.
$excel.ActiveWindow.ScrollRow = 123
$ws.Range("A141").Select()
